$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "29.778.88"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -0.49%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.862.12"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -1.57%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.004"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  +0.28%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "0.7327"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -5.28%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "241.24"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.98%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "1.004"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +0.33%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3081"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -1.69%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "24.44"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -4.66%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.07034"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -4.17%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.08399"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +4.21%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.7445"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -3.52%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "1.872.77"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +0.21%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "5.309"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -3.38%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "91.92"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -2.25%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "29.794.34"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -0.39%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "6.065"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -2.35%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "13.47"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -3.76%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "239.19"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -2.91%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "0.000007752"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -1.32%  "

$ws.Range("E21").Value = "  +0.21%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "2.138.05"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +0.32%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "1.004"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +0.30%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "7.879"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -3.25%  "

$ws.Range("E25").Value = "  -0.47%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "9.235"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -2.12%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "161.95"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -0.10%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "18.45"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -1.59%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "1.993"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -1.48%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "1.487"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +4.54%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "1.530"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -0.68%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "4.422"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -1.01%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "4.111"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +1.24%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.05336"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -3.93%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.223"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -0.88%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.7393"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -1.35%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +0.02%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "2.699"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +0.61%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.01918"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -0.58%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "2.739"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -1.81%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.4401"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -1.41%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "1.098.57"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -0.17%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "5.997"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -0.15%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "71.48"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -3.78%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.8594"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +1.07%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "101.83"
$cell.Style = "Normal"

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "7.648"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +1.58%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "1.819"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -3.53%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "2.986"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +0.02%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "2.038.43"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +0.24%  "
